# Sync attendance_reports: move the "System" entry from the front to the
# back of the "Recorded By" (column G) comma-separated list.
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# "Recorded By" is column G (7) per the header row, but resolve it
# dynamically in case column order ever changes.
$recordedByCol = 7
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Text -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text

    if ($val -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value = "$rest, System"
    }
}
